$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param(
        $ws,
        [int]$row,
        [double[]]$vals
    )
    $cols = @("A","B","C","D","E")
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

# --- Sheet: Escapement ---
$ws = $wb.Worksheets.Item("Escapement")
Set-RowValues $ws 59 @(4028889.48167101, 2385881.40802501, 4110477.96951001, 2412280.997283, 1732176.2079256)

# --- Sheet: Total Catch ---
$ws = $wb.Worksheets.Item("Total Catch")
Set-RowValues $ws 57 @(3067971.65104639, 1292123.43158704, 4704909.05283045, 4495955.34878523, 3602430.96056782)
Set-RowValues $ws 58 @(5565069.48341984, 1034687.52937156, 5913187.20470384, 13093842.4110679, 2057215.09245664)
Set-RowValues $ws 59 @(5944787.03766096, 2227286.75334413, 6101073.0695064, 13404240.0411759, 3396471.7587296)

# --- Sheet: Run Size ---
$ws = $wb.Worksheets.Item("Run Size")
Set-RowValues $ws 57 @(7466679.6510669, 2873549.43197125, 6926061.05364625, 6104312.34861969, 4770222.96091154)
Set-RowValues $ws 58 @(7936311.48368038, 1855145.52936893, 8824657.20550253, 15434052.4202096, 3604963.09288495)
Set-RowValues $ws 59 @(9975755.03733196, 4613804.75336914, 10213233.0700164, 15793968.0414589, 5142411.7586552)

# --- Sheet: Run Size no Offshore ---
$ws = $wb.Worksheets.Item("Run Size no Offshore")
Set-RowValues $ws 58 @(7857381.16684172, 1836550.1266757, 8737379.42620984, 15251550.7806298, 3562468.037446)
Set-RowValues $ws 59 @(9811353.21886345, 4534935.15862568, 10046645.3940871, 15501887.0430196, 5047701.89071117)
